# Update gh-pages to output generated at 456a3b4
#
# Sheet "展览" (sheet index 1) and sheet "全部类型" (sheet index 4) both get:
#   - five small "want-to-go count" (column F) bumps on existing rows
#   - a brand-new row inserted for 新余·2024第三届MG动漫嘉年华 (2024-07-12),
#     pushing every later row down by one
#   - the want-to-go count for 赣州·第二届异次元动漫嘉年华 (now the last row)
#     bumped from 444 to 447

$wb = $excel.ActiveWorkbook

function Update-ExhibitionSheet($ws, $newRow) {

    $lastRow = $newRow + 10   # 赣州·第二届异次元动漫嘉年华 ends up 10 rows below the insert

    # --- small numeric bumps on rows above the insertion point ---------------
    $ws.Cells.Item($newRow - 17, 6).Value = 179    # 九江·首届萤火之星国风动漫嘉年华: 178 -> 179
    $ws.Cells.Item($newRow - 16, 6).Value = 2720   # 南昌·CM02动漫游戏博览会: 2708 -> 2720
    $ws.Cells.Item($newRow - 14, 6).Value = 96     # 九江·第四届ACD动漫游戏嘉年华: 95 -> 96
    $ws.Cells.Item($newRow - 11, 6).Value = 552    # 南昌·第三届龙年动漫展: 551 -> 552
    $ws.Cells.Item($newRow - 2,  6).Value = 218    # 宜春·BM次元盛典运动番only: 217 -> 218

    # --- insert the new row, carrying down the formatting of the row above ---
    $ws.Rows.Item($newRow).Insert()
    $ws.Range($ws.Cells.Item($newRow - 1, 1), $ws.Cells.Item($newRow - 1, 9)).Copy()
    $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 9)).PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    $ws.Cells.Item($newRow, 1).Value = $newRow - 1
    # Column B holds plain date-looking text ("2024-07-12") elsewhere in the
    # sheet; force text format first so Excel doesn't silently convert it to
    # a date serial number.
    $ws.Cells.Item($newRow, 2).NumberFormat = "@"
    $ws.Cells.Item($newRow, 2).Value = "2024-07-12"
    $ws.Cells.Item($newRow, 3).Value = "新余·2024第三届MG动漫嘉年华"
    $ws.Cells.Item($newRow, 4).Value = "仙女湖大道与五一南路交叉口西约180米 老上海风情街水晶厅"
    $ws.Cells.Item($newRow, 5).Value = "2024.07.12 10:00-07.13 17:30"
    $ws.Cells.Item($newRow, 6).Value = 1
    $ws.Cells.Item($newRow, 7).Value = 55
    $ws.Cells.Item($newRow, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86536"
    $ws.Cells.Item($newRow, 9).Value = "//i0.hdslb.com/bfs/openplatform/202405/11RbfeFq1716813676323.jpeg"

    # --- 赣州·第二届异次元动漫嘉年华, now pushed down to $lastRow: 444 -> 447 --
    $ws.Cells.Item($lastRow, 6).Value = 447
}

# Sheet 1: 展览 — new row lands at row 22
Update-ExhibitionSheet $wb.Worksheets.Item(1) 22

# Sheet 4: 全部类型 — same edits, shifted one row down (extra row at top)
Update-ExhibitionSheet $wb.Worksheets.Item(4) 23
